$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Create the hyperlink for the new problem name first (TextToDisplay only
# controls the stored <hyperlink display="..."/> attribute - it also
# clobbers the cell's text, so the real label is written back afterwards).
$ws.Hyperlinks.Add($ws.Range("B37"), "https://leetcode.com/problems/majority-element/", [System.Type]::Missing, [System.Type]::Missing, "https://leetcode.com/problems/majority-element/") | Out-Null

# New row of LeetCode stats data: "Majority Element"
$ws.Range("B37").Value = "Majority Element"
$ws.Range("B37").Style = $ws.Range("B36").Style
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 136
$ws.Range("F37").Value = 0.0047
$ws.Range("G37").Value = 17.8
$ws.Range("H37").Value = 0.1742
$ws.Range("I37").Value = "https://leetcode.com/problems/majority-element/submissions/1075778041/"

# Mirror the saved selection state on Sheet1
$ws.Range("I42").Select()
